# Update cryptos list values - GitHub Actions refresh (Thu Jul 11 18:31:32 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.817.56"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.138.75"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.54"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.46"
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.137.82"
$ws.Range("E8").Value = "  +1.28%  "
$ws.Range("E9").Value = "  +2.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.22"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.398"
$ws.Range("E12").Value = "  +3.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.677.77"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("E14").Value = "  +2.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.58"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000166"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.969.01"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.142.29"
$ws.Range("E19").Value = "  -2.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.77"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.95"
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "353.18"
$ws.Range("E22").Value = "  +5.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.66"
$ws.Range("E24").Value = "  +2.83%  "
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0918"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.54"
$ws.Range("E29").Value = "  +4.69%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.19"
$ws.Range("E31").Value = "  -5.20%  "
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.20"
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.04"
$ws.Range("E34").Value = "  +8.48%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.18"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.22"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.17"
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.43"
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.28"
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0672"
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.21"
$ws.Range("E41").Value = "  +7.09%  "
$ws.Range("E42").Value = "  +6.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.704"
$ws.Range("E43").Value = "  +2.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.179.77"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.58"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0271"
$ws.Range("E46").Value = "  +4.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.332.98"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.969"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.42"
$ws.Range("E51").Value = "  -1.68%  "
